$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.493.21"
$ws.Range("E2").Value = "'  +0.31%  "
$ws.Range("D3").Value = "'1.913.93"
$ws.Range("E3").Value = "'  -0.10%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'244.60"
$ws.Range("E5").Value = "'  +1.52%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("D7").Value = "'0.4836"
$ws.Range("E7").Value = "'  +3.30%  "
$ws.Range("D8").Value = "'0.2898"
$ws.Range("E8").Value = "'  +1.89%  "
$ws.Range("D9").Value = "'0.06726"
$ws.Range("E9").Value = "'  -0.99%  "
$ws.Range("D10").Value = "'109.89"
$ws.Range("D11").Value = "'19.00"
$ws.Range("E11").Value = "'  +4.85%  "
$ws.Range("D12").Value = "'1.918.58"
$ws.Range("E12").Value = "'  +0.22%  "
$ws.Range("D13").Value = "'0.07549"
$ws.Range("E13").Value = "'  -1.00%  "
$ws.Range("D14").Value = "'5.279"
$ws.Range("E14").Value = "'  +2.00%  "
$ws.Range("D15").Value = "'0.6728"
$ws.Range("E15").Value = "'  +2.91%  "
$ws.Range("D16").Value = "'282.38"
$ws.Range("E16").Value = "'  -1.93%  "
$ws.Range("D17").Value = "'30.509.36"
$ws.Range("E17").Value = "'  +0.36%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "'  +0.05%  "
$ws.Range("D19").Value = "'0.000007568"
$ws.Range("E19").Value = "'  -0.56%  "
$ws.Range("E20").Value = "'  -0.66%  "
$ws.Range("D21").Value = "'5.515"
$ws.Range("E21").Value = "'  +5.77%  "
$ws.Range("D22").Value = "'2.166.79"
$ws.Range("E22").Value = "'  +0.51%  "
$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "'  -0.15%  "
$ws.Range("D24").Value = "'6.471"
$ws.Range("E24").Value = "'  +4.61%  "
$ws.Range("D25").Value = "'9.468"
$ws.Range("E25").Value = "'  +2.40%  "
$ws.Range("D26").Value = "'164.31"
$ws.Range("E26").Value = "'  -2.23%  "
$ws.Range("D27").Value = "'20.26"
$ws.Range("E27").Value = "'  -6.55%  "
$ws.Range("D28").Value = "'2.127"
$ws.Range("E28").Value = "'  +4.33%  "
$ws.Range("D29").Value = "'0.1055"
$ws.Range("E29").Value = "'  -1.11%  "
$ws.Range("D30").Value = "'1.400"
$ws.Range("E30").Value = "'  +2.18%  "
$ws.Range("E31").Value = "'  +0.36%  "
$ws.Range("D32").Value = "'4.045"
$ws.Range("E32").Value = "'  +2.82%  "
$ws.Range("D33").Value = "'0.04998"
$ws.Range("E33").Value = "'  -0.53%  "
$ws.Range("D34").Value = "'0.7313"
$ws.Range("E34").Value = "'  -0.62%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "'  -0.76%  "
$ws.Range("D36").Value = "'0.9994"
$ws.Range("E36").Value = "'  +0.05%  "
$ws.Range("D37").Value = "'2.728"
$ws.Range("E37").Value = "'  -0.33%  "
$ws.Range("E38").Value = "'  +0.02%  "
$ws.Range("E39").Value = "'  -0.75%  "
$ws.Range("D40").Value = "'110.92"
$ws.Range("E40").Value = "'  +2.11%  "
$ws.Range("E41").Value = "'  -1.60%  "
$ws.Range("D42").Value = "'0.4449"
$ws.Range("E42").Value = "'  +5.95%  "
$ws.Range("D43").Value = "'0.8663"
$ws.Range("E43").Value = "'  -0.75%  "
$ws.Range("D44").Value = "'5.803"
$ws.Range("E44").Value = "'  -0.57%  "
$ws.Range("D45").Value = "'0.9994"
$ws.Range("E45").Value = "'  +0.00%  "
$ws.Range("D46").Value = "'68.13"
$ws.Range("E46").Value = "'  +0.94%  "
$ws.Range("E47").Value = "'  +2.78%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'48.85"
$ws.Range("E48").Value = "'  -7.46%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.287"
$ws.Range("E49").Value = "'  +0.78%  "
$ws.Range("D50").Value = "'0.1239"
$ws.Range("E50").Value = "'  +2.78%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'34.79"
$ws.Range("E51").Value = "'  +0.42%  "
